# Updated capital structure database
# Applies new values to the data rows (rows 2 and 3) of the active sheet,
# matching the updated figures in the source database export.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "G"  = 0.2008165348303139
    "H"  = 0.2008165348303139
    "I"  = 0.1605001275835672
    "J"  = 0.1459953023889033
    "K"  = 58.9
    "L"  = 0.1502934422046441
    "M"  = 81.2
    "N"  = 0.1075924208294687
    "O"  = 1.378607809847199
    "P"  = 81.2
    "Q"  = 0.1075924208294687
    "R"  = 1.378607809847199
    "S"  = 0
    "T"  = 0
    "U"  = 38.6
    "V"  = 0.05114615078839273
    "W"  = 0.4153737658674189
    "X"  = 0.05889345716520324
    "Y"  = 0.3564803087022156
    "Z"  = 4.186518534344621
    "AA" = 0.6112120393783909
    "AB" = 0.05862708088799806
    "AC" = 0.5525849584903928
    "AD" = 5.25
    "AE" = 0
    "AF" = 5.25
    "AG" = -33.35
    "AH" = 0.006908349233502204
    "AI" = 0.03550896178559351
    "AJ" = -0.04623275802315104
    "AK" = -0.3052631578947368
    "AL" = 0
    "AM" = -1.66
    "AN" = 0.07894736842105263
    "AP" = -0.5015037593984962
    "AQ" = -37.89156626506024
}

foreach ($row in @(2, 3)) {
    foreach ($col in $values.Keys) {
        $ws.Range("$col$row").Value = $values[$col]
    }
}
